# Fruta / hortaliza, semanal
# Update weekly price records (rows 11-15) and append two additional
# records that were missing (now rows 16-17) on the Chirimoya sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: date, quality, volume & price updates ---------------------
$ws.Range("D11").Value = 44874
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 30
$ws.Range("N11").Value = 25000
$ws.Range("O11").Value = 25000
$ws.Range("P11").Value = 25000
$ws.Range("S11").Value = 2500

# --- Row 12: date, quality, volume & price updates ---------------------
$ws.Range("D12").Value = 44874
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 80
$ws.Range("N12").Value = 23000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 23500
$ws.Range("S12").Value = 2350

# --- Row 13: date & volume updates --------------------------------------
$ws.Range("D13").Value = 44448
$ws.Range("M13").Value = 60

# --- Row 14: date, quality & price updates ------------------------------
$ws.Range("D14").Value = 44848
$ws.Range("L14").Value = "Especial"
$ws.Range("N14").Value = 24000
$ws.Range("O14").Value = 25000
$ws.Range("P14").Value = 24500
$ws.Range("S14").Value = 2450

# --- Row 15: date, volume & price updates -------------------------------
$ws.Range("D15").Value = 44848
$ws.Range("M15").Value = 120
$ws.Range("N15").Value = 21000
$ws.Range("O15").Value = 22000
$ws.Range("P15").Value = 21500
$ws.Range("S15").Value = 2150

# --- Row 16 (new record) -------------------------------------------------
$ws.Range("A16").Value = 7
$ws.Range("B16").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C16").Value = "Ñuble"
$ws.Range("D16").Value = 44452
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E16").Value = 16
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100107
$ws.Range("H16").Value = "Otros"
$ws.Range("I16").Value = 100107002
$ws.Range("J16").Value = "Chirimoya"
$ws.Range("K16").Value = "Cultivar IV Región"
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 21000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 21500
$ws.Range("Q16").Value = "$/bandeja 10 kilos"
$ws.Range("R16").Value = "Provincia de Limarí"
$ws.Range("S16").Value = 2150
$ws.Range("T16").Value = 10

# --- Row 17 (new record) -------------------------------------------------
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44487
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100107
$ws.Range("H17").Value = "Otros"
$ws.Range("I17").Value = 100107002
$ws.Range("J17").Value = "Chirimoya"
$ws.Range("K17").Value = "Cultivar IV Región"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 30
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 23500
$ws.Range("Q17").Value = "$/bandeja 10 kilos"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 2350
$ws.Range("T17").Value = 10
